$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$n = 182
$arr = New-Object 'object[,]' $n,3
$arr[0,0] = 'scientific.name'; $arr[0,1] = 'status'; $arr[0,2] = 'in.anp'
$arr[1,0] = 'Actaea pachypoda'; $arr[1,1] = 'rare native'; $arr[1,2] = 'P'
$arr[2,0] = 'Actaea rubra'; $arr[2,1] = 'rare native'; $arr[2,2] = 'P'
$arr[3,0] = 'Agalinis purpurea'; $arr[3,1] = 'rare native'; $arr[3,2] = 'P'
$arr[4,0] = 'Amphicarpaea bracteata'; $arr[4,1] = 'rare native'; $arr[4,2] = 'P'
$arr[5,0] = 'Andropogon gerardi'; $arr[5,1] = 'rare native'; $arr[5,2] = 'P'
$arr[6,0] = 'Antennaria neglecta'; $arr[6,1] = 'rare native'; $arr[6,2] = 'P'
$arr[7,0] = 'Arisaema triphyllum'; $arr[7,1] = 'rare native'; $arr[7,2] = 'P'
$arr[8,0] = 'Asclepias syriaca'; $arr[8,1] = 'rare native'; $arr[8,2] = 'P'
$arr[9,0] = 'Asplenium trichomanes'; $arr[9,1] = 'rare native'; $arr[9,2] = 'P'
$arr[10,0] = 'Atriplex acadiensis'; $arr[10,1] = 'rare native'; $arr[10,2] = 'P'
$arr[11,0] = 'Atriplex glabriuscula'; $arr[11,1] = 'rare native'; $arr[11,2] = 'P'
$arr[12,0] = 'Atriplex dioica'; $arr[12,1] = 'rare native'; $arr[12,2] = 'P'
$arr[13,0] = 'Bartonia paniculata'; $arr[13,1] = 'rare native'; $arr[13,2] = 'P'
$arr[14,0] = 'Bartonia virginica'; $arr[14,1] = 'rare native'; $arr[14,2] = 'P'
$arr[15,0] = 'Calamagrostis inexpansa'; $arr[15,1] = 'rare native'; $arr[15,2] = 'P'
$arr[16,0] = 'Calamagrostis stricta'; $arr[16,1] = 'rare native'; $arr[16,2] = 'P'
$arr[17,0] = 'Calluna vulgaris'; $arr[17,1] = 'rare native'; $arr[17,2] = 'P'
$arr[18,0] = 'Caltha palustris'; $arr[18,1] = 'rare native'; $arr[18,2] = 'P'
$arr[19,0] = 'Capnoides sempervirens'; $arr[19,1] = 'rare native'; $arr[19,2] = 'P'
$arr[20,0] = 'Carex houghtoniana'; $arr[20,1] = 'rare native'; $arr[20,2] = 'P'
$arr[21,0] = 'Carex ormostachya'; $arr[21,1] = 'rare native'; $arr[21,2] = 'P'
$arr[22,0] = 'Carex wiegandii'; $arr[22,1] = 'rare native'; $arr[22,2] = 'P'
$arr[23,0] = 'Chelone glabra'; $arr[23,1] = 'rare native'; $arr[23,2] = 'P'
$arr[24,0] = 'Cicuta maculata'; $arr[24,1] = 'rare native'; $arr[24,2] = 'P'
$arr[25,0] = 'Corallorhiza maculata'; $arr[25,1] = 'rare native'; $arr[25,2] = 'P'
$arr[26,0] = 'Corallorhiza trifida'; $arr[26,1] = 'rare native'; $arr[26,2] = 'P'
$arr[27,0] = 'Corema conradii'; $arr[27,1] = 'rare native'; $arr[27,2] = 'P'
$arr[28,0] = 'Crataegus macracantha'; $arr[28,1] = 'rare native'; $arr[28,2] = 'P'
$arr[29,0] = 'Cuscuta gronovii'; $arr[29,1] = 'rare native'; $arr[29,2] = 'P'
$arr[30,0] = 'Cypripedium acaule'; $arr[30,1] = 'rare native'; $arr[30,2] = 'P'
$arr[31,0] = 'Cystopteris fragilis'; $arr[31,1] = 'rare native'; $arr[31,2] = 'P'
$arr[32,0] = 'Dasiphora fruticosa'; $arr[32,1] = 'rare native'; $arr[32,2] = 'P'
$arr[33,0] = 'Diphasiastrum digitatum'; $arr[33,1] = 'rare native'; $arr[33,2] = 'P'
$arr[34,0] = 'Epigaea repens'; $arr[34,1] = 'rare native'; $arr[34,2] = 'P'
$arr[35,0] = 'Epilobium ciliatum'; $arr[35,1] = 'rare native'; $arr[35,2] = 'P'
$arr[36,0] = 'Epilobium coloratum'; $arr[36,1] = 'rare native'; $arr[36,2] = 'P'
$arr[37,0] = 'Epilobium palustre'; $arr[37,1] = 'rare native'; $arr[37,2] = 'P'
$arr[38,0] = 'Equisetum fluviatile'; $arr[38,1] = 'rare native'; $arr[38,2] = 'P'
$arr[39,0] = 'Euthamia graminifolia'; $arr[39,1] = 'rare native'; $arr[39,2] = 'P'
$arr[40,0] = 'Eutrochium maculatum'; $arr[40,1] = 'rare native'; $arr[40,2] = 'P'
$arr[41,0] = 'Geum rivale'; $arr[41,1] = 'rare native'; $arr[41,2] = 'P'
$arr[42,0] = 'Goodyera pubescens'; $arr[42,1] = 'rare native'; $arr[42,2] = 'P'
$arr[43,0] = 'Hieracium scabrum'; $arr[43,1] = 'rare native'; $arr[43,2] = 'P'
$arr[44,0] = 'Huperzia appressa'; $arr[44,1] = 'rare native'; $arr[44,2] = 'P'
$arr[45,0] = 'Monotropa hypopitys'; $arr[45,1] = 'rare native'; $arr[45,2] = 'P'
$arr[46,0] = 'Lilium philadelphicum'; $arr[46,1] = 'rare native'; $arr[46,2] = 'P'
$arr[47,0] = 'Lonicera canadensis'; $arr[47,1] = 'rare native'; $arr[47,2] = 'P'
$arr[48,0] = 'Lycopodiella inundata'; $arr[48,1] = 'rare native'; $arr[48,2] = 'P'
$arr[49,0] = 'Lysimachia maritima'; $arr[49,1] = 'rare native'; $arr[49,2] = 'P'
$arr[50,0] = 'Maianthemum racemosum'; $arr[50,1] = 'rare native'; $arr[50,2] = 'P'
$arr[51,0] = 'Medeola virginiana'; $arr[51,1] = 'rare native'; $arr[51,2] = 'P'
$arr[52,0] = 'Mertensia maritima'; $arr[52,1] = 'rare native'; $arr[52,2] = 'P'
$arr[53,0] = 'Mollugo verticillata'; $arr[53,1] = 'rare native'; $arr[53,2] = 'P'
$arr[54,0] = 'Panicum virgatum'; $arr[54,1] = 'rare native'; $arr[54,2] = 'P'
$arr[55,0] = 'Pinus resinosa'; $arr[55,1] = 'rare native'; $arr[55,2] = 'P'
$arr[56,0] = 'Platanthera clavellata'; $arr[56,1] = 'rare native'; $arr[56,2] = 'P'
$arr[57,0] = 'Polygala sanguinea'; $arr[57,1] = 'rare native'; $arr[57,2] = 'P'
$arr[58,0] = 'Polygonatum pubescens'; $arr[58,1] = 'rare native'; $arr[58,2] = 'P'
$arr[59,0] = 'Polypodium appalachianum'; $arr[59,1] = 'rare native'; $arr[59,2] = 'P'
$arr[60,0] = 'Proserpinaca pectinata'; $arr[60,1] = 'rare native'; $arr[60,2] = 'P'
$arr[61,0] = 'Ranunculus recurvatus'; $arr[61,1] = 'rare native'; $arr[61,2] = 'P'
$arr[62,0] = 'Rhexia virginica'; $arr[62,1] = 'rare native'; $arr[62,2] = 'P'
$arr[63,0] = 'Rhodiola rosea'; $arr[63,1] = 'rare native'; $arr[63,2] = 'P'
$arr[64,0] = 'Rhynchospora capitellata'; $arr[64,1] = 'rare native'; $arr[64,2] = 'P'
$arr[65,0] = 'Rubus chamaemorus'; $arr[65,1] = 'rare native'; $arr[65,2] = 'P'
$arr[66,0] = 'Rubus dalibarda'; $arr[66,1] = 'rare native'; $arr[66,2] = 'P'
$arr[67,0] = 'Solidago flexicaulis'; $arr[67,1] = 'rare native'; $arr[67,2] = 'P'
$arr[68,0] = 'Sorghastrum nutans'; $arr[68,1] = 'rare native'; $arr[68,2] = 'P'
$arr[69,0] = 'Sparganium emersum'; $arr[69,1] = 'rare native'; $arr[69,2] = 'P'
$arr[70,0] = 'Spergularia canadensis'; $arr[70,1] = 'rare native'; $arr[70,2] = 'P'
$arr[71,0] = 'Spergularia marina'; $arr[71,1] = 'rare native'; $arr[71,2] = 'P'
$arr[72,0] = 'Spinulum annotinum'; $arr[72,1] = 'rare native'; $arr[72,2] = 'P'
$arr[73,0] = 'Spiranthes cernua'; $arr[73,1] = 'rare native'; $arr[73,2] = 'P'
$arr[74,0] = 'Suaeda calceoliformis'; $arr[74,1] = 'rare native'; $arr[74,2] = 'P'
$arr[75,0] = 'Cornus rugosa'; $arr[75,1] = 'rare native'; $arr[75,2] = 'P'
$arr[76,0] = 'Symplocarpus foetidus'; $arr[76,1] = 'rare native'; $arr[76,2] = 'P'
$arr[77,0] = 'Taxus canadensis'; $arr[77,1] = 'rare native'; $arr[77,2] = 'P'
$arr[78,0] = 'Teucrium canadense'; $arr[78,1] = 'rare native'; $arr[78,2] = 'P'
$arr[79,0] = 'Toxicodendron radicans'; $arr[79,1] = 'rare native'; $arr[79,2] = 'P'
$arr[80,0] = 'Toxicodendron rydbergii'; $arr[80,1] = 'rare native'; $arr[80,2] = 'P'
$arr[81,0] = 'Trichophorum cespitosum'; $arr[81,1] = 'rare native'; $arr[81,2] = 'P'
$arr[82,0] = 'Trillium cernuum'; $arr[82,1] = 'rare native'; $arr[82,2] = 'P'
$arr[83,0] = 'Trillium undulatum'; $arr[83,1] = 'rare native'; $arr[83,2] = 'P'
$arr[84,0] = 'Typha latifolia'; $arr[84,1] = 'rare native'; $arr[84,2] = 'P'
$arr[85,0] = 'Uvularia sessilifolia'; $arr[85,1] = 'rare native'; $arr[85,2] = 'P'
$arr[86,0] = 'Verbena hastata'; $arr[86,1] = 'rare native'; $arr[86,2] = 'P'
$arr[87,0] = 'Viburnum dentatum'; $arr[87,1] = 'rare native'; $arr[87,2] = 'P'
$arr[88,0] = 'Viburnum recognitum'; $arr[88,1] = 'rare native'; $arr[88,2] = 'P'
$arr[89,0] = 'Woodsia ilvensis'; $arr[89,1] = 'rare native'; $arr[89,2] = 'P'
$arr[90,0] = 'Woodwardia virginica'; $arr[90,1] = 'rare native'; $arr[90,2] = 'P'
$arr[91,0] = 'Xyris difformis'; $arr[91,1] = 'rare native'; $arr[91,2] = 'P'
$arr[92,0] = 'Phellodendron amurense'; $arr[92,1] = 'invasive not established'; $arr[92,2] = 'Y'
$arr[93,0] = 'Anthriscus sylvestris'; $arr[93,1] = 'invasive not established'; $arr[93,2] = 'N'
$arr[94,0] = 'Vincetoxicum nigrum'; $arr[94,1] = 'invasive not established'; $arr[94,2] = 'N'
$arr[95,0] = 'Vincetoxicum rossicum'; $arr[95,1] = 'invasive not established'; $arr[95,2] = 'N'
$arr[96,0] = 'Dioscorea polystachya'; $arr[96,1] = 'invasive not established'; $arr[96,2] = 'N'
$arr[97,0] = 'Akebia quinata'; $arr[97,1] = 'invasive not established'; $arr[97,2] = 'N'
$arr[98,0] = 'Alnus glutinosa'; $arr[98,1] = 'invasive not established'; $arr[98,2] = 'N'
$arr[99,0] = 'Amorpha fruticosa'; $arr[99,1] = 'invasive not established'; $arr[99,2] = 'N'
$arr[100,0] = 'Butomus umbellatus'; $arr[100,1] = 'invasive not established'; $arr[100,2] = 'N'
$arr[101,0] = 'Augopodium podagraria'; $arr[101,1] = 'invasive not established'; $arr[101,2] = 'N'
$arr[102,0] = 'Microstegeum vimineum'; $arr[102,1] = 'invasive not established'; $arr[102,2] = 'N'
$arr[103,0] = 'Lonicera tatarica'; $arr[103,1] = 'invasive not established'; $arr[103,2] = 'N'
$arr[104,0] = 'Lonicera periclymenum'; $arr[104,1] = 'invasive not established'; $arr[104,2] = 'Y'
$arr[105,0] = 'Lonicera maackii'; $arr[105,1] = 'invasive not established'; $arr[105,2] = 'N'
$arr[106,0] = 'Persicaria perfoliata'; $arr[106,1] = 'invasive not established'; $arr[106,2] = 'N'
$arr[107,0] = 'Berberis vulgaris'; $arr[107,1] = 'invasive not established'; $arr[107,2] = 'N'
$arr[108,0] = 'Impatiens glandulifera'; $arr[108,1] = 'invasive not established'; $arr[108,2] = 'N'
$arr[109,0] = 'Lepidium latifolium'; $arr[109,1] = 'invasive not established'; $arr[109,2] = 'N'
$arr[110,0] = 'Phragmites australis'; $arr[110,1] = 'invasive not established'; $arr[110,2] = 'N'
$arr[111,0] = 'Ampelopsis glandulosa'; $arr[111,1] = 'invasive not established'; $arr[111,2] = 'N'
$arr[112,0] = 'Rhamnus cathartica'; $arr[112,1] = 'invasive not established'; $arr[112,2] = 'N'
$arr[113,0] = 'Ailanthus altissima'; $arr[113,1] = 'invasive not established'; $arr[113,2] = 'N'
$arr[114,0] = 'Typha angustifolia'; $arr[114,1] = 'invasive not established'; $arr[114,2] = 'N'
$arr[115,0] = 'Oplismenus hirtellus'; $arr[115,1] = 'invasive not established'; $arr[115,2] = 'N'
$arr[116,0] = 'Populus alba'; $arr[116,1] = 'invasive not established'; $arr[116,2] = 'Y'
$arr[117,0] = 'Rubus phoenicolasius'; $arr[117,1] = 'invasive not established'; $arr[117,2] = 'N'
$arr[118,0] = 'Acer ginnala'; $arr[118,1] = 'invasive established'; $arr[118,2] = 'Y'
$arr[119,0] = 'Acer platanoides'; $arr[119,1] = 'invasive established'; $arr[119,2] = 'Y'
$arr[120,0] = 'Achillea ptarmica'; $arr[120,1] = 'invasive established'; $arr[120,2] = 'Y'
$arr[121,0] = 'Alliaria petiolata'; $arr[121,1] = 'invasive established'; $arr[121,2] = 'Y'
$arr[122,0] = 'Berberis thunbergii'; $arr[122,1] = 'invasive established'; $arr[122,2] = 'Y'
$arr[123,0] = 'Cardamine impatiens'; $arr[123,1] = 'invasive established'; $arr[123,2] = 'Y'
$arr[124,0] = 'Celastrus orbiculatus'; $arr[124,1] = 'invasive established'; $arr[124,2] = 'Y'
$arr[125,0] = 'Centaurea stoebe'; $arr[125,1] = 'invasive established'; $arr[125,2] = 'Y'
$arr[126,0] = 'Centaurea nigrescens'; $arr[126,1] = 'invasive established'; $arr[126,2] = 'Y'
$arr[127,0] = 'Cirsium arvense'; $arr[127,1] = 'invasive established'; $arr[127,2] = 'Y'
$arr[128,0] = 'Cirsium vulgare'; $arr[128,1] = 'invasive established'; $arr[128,2] = 'Y'
$arr[129,0] = 'Calluna vulgaris'; $arr[129,1] = 'invasive established'; $arr[129,2] = 'Y'
$arr[130,0] = 'Digitalis purpurea'; $arr[130,1] = 'invasive established'; $arr[130,2] = 'Y'
$arr[131,0] = 'Elaeagnus umbellata'; $arr[131,1] = 'invasive established'; $arr[131,2] = 'Y'
$arr[132,0] = 'Epipactis helleborine'; $arr[132,1] = 'invasive established'; $arr[132,2] = 'Y'
$arr[133,0] = 'Euonymus alatus'; $arr[133,1] = 'invasive established'; $arr[133,2] = 'Y'
$arr[134,0] = 'Euphorbia cyparissias'; $arr[134,1] = 'invasive established'; $arr[134,2] = 'Y'
$arr[135,0] = 'Galeopsis bifida'; $arr[135,1] = 'invasive established'; $arr[135,2] = 'Y'
$arr[136,0] = 'Genista tinctoria'; $arr[136,1] = 'invasive established'; $arr[136,2] = 'Y'
$arr[137,0] = 'Hypericum prolificum'; $arr[137,1] = 'invasive established'; $arr[137,2] = 'Y'
$arr[138,0] = 'Hylotelephium telephium'; $arr[138,1] = 'invasive established'; $arr[138,2] = 'Y'
$arr[139,0] = 'Lotus corniculatus'; $arr[139,1] = 'invasive established'; $arr[139,2] = 'Y'
$arr[140,0] = 'Reynoutria japonica'; $arr[140,1] = 'invasive established'; $arr[140,2] = 'Y'
$arr[141,0] = 'Heracleum mantegazzianum'; $arr[141,1] = 'invasive established'; $arr[141,2] = 'Y'
$arr[142,0] = 'Iris pseudacorus'; $arr[142,1] = 'invasive established'; $arr[142,2] = 'Y'
$arr[143,0] = 'Impatiens glandulifera'; $arr[143,1] = 'invasive established'; $arr[143,2] = 'Y'
$arr[144,0] = 'Jacobaea vulgaris'; $arr[144,1] = 'invasive established'; $arr[144,2] = 'Y'
$arr[145,0] = 'Ligustrum spp.'; $arr[145,1] = 'invasive established'; $arr[145,2] = 'Y'
$arr[146,0] = 'Lonicera japonica'; $arr[146,1] = 'invasive established'; $arr[146,2] = 'Y'
$arr[147,0] = 'Lupinus polyphyllus'; $arr[147,1] = 'invasive established'; $arr[147,2] = 'Y'
$arr[148,0] = 'Lythrum salicaria'; $arr[148,1] = 'invasive established'; $arr[148,2] = 'Y'
$arr[149,0] = 'Frangula alnus'; $arr[149,1] = 'invasive established'; $arr[149,2] = 'Y'
$arr[150,0] = 'Potentilla recta'; $arr[150,1] = 'invasive established'; $arr[150,2] = 'Y'
$arr[151,0] = 'Rosa multiflora'; $arr[151,1] = 'invasive established'; $arr[151,2] = 'Y'
$arr[152,0] = 'Rosa rugosa'; $arr[152,1] = 'invasive established'; $arr[152,2] = 'Y'
$arr[153,0] = 'Securigera varia'; $arr[153,1] = 'invasive established'; $arr[153,2] = 'Y'
$arr[154,0] = 'Sisymbrium officinale'; $arr[154,1] = 'invasive established'; $arr[154,2] = 'Y'
$arr[155,0] = 'Solanum dulcamara'; $arr[155,1] = 'invasive established'; $arr[155,2] = 'Y'
$arr[156,0] = 'Spergularia rubra'; $arr[156,1] = 'invasive established'; $arr[156,2] = 'Y'
$arr[157,0] = 'Trifolium hybridum'; $arr[157,1] = 'invasive established'; $arr[157,2] = 'Y'
$arr[158,0] = 'Tussilago farfara'; $arr[158,1] = 'invasive established'; $arr[158,2] = 'Y'
$arr[159,0] = 'Verbascum thapsus'; $arr[159,1] = 'invasive established'; $arr[159,2] = 'Y'
$arr[160,0] = 'Adelges tsugae'; $arr[160,1] = 'pest disease'; $arr[160,2] = 'Y'
$arr[161,0] = 'Adelges piceae'; $arr[161,1] = 'pest disease'; $arr[161,2] = 'Y'
$arr[162,0] = 'Pyrrhalta viburni'; $arr[162,1] = 'pest disease'; $arr[162,2] = 'Y'
$arr[163,0] = 'Litylenchus crenatae'; $arr[163,1] = 'pest disease'; $arr[163,2] = 'Y'
$arr[164,0] = 'Anoplophora glabripennis'; $arr[164,1] = 'pest disease'; $arr[164,2] = 'N'
$arr[165,0] = 'Agrilus planipennis'; $arr[165,1] = 'pest disease'; $arr[165,2] = 'N'
$arr[166,0] = 'Euproctis chrysorrhoea'; $arr[166,1] = 'pest disease'; $arr[166,2] = 'Y'
$arr[167,0] = 'Operophtera brumata'; $arr[167,1] = 'pest disease'; $arr[167,2] = 'Y'
$arr[168,0] = 'Bretziella fagacearum'; $arr[168,1] = 'pest disease'; $arr[168,2] = 'N'
$arr[169,0] = 'Choristoneura fumiferana'; $arr[169,1] = 'pest disease'; $arr[169,2] = 'N'
$arr[170,0] = 'Tetropium castaneum'; $arr[170,1] = 'pest disease'; $arr[170,2] = 'N'
$arr[171,0] = 'Fiorinia externa'; $arr[171,1] = 'pest disease'; $arr[171,2] = 'N'
$arr[172,0] = 'Matsucoccus matsumurae'; $arr[172,1] = 'pest disease'; $arr[172,2] = 'Y'
$arr[173,0] = 'Amynthas agrestis'; $arr[173,1] = 'pest disease'; $arr[173,2] = 'Y'
$arr[174,0] = 'Amynthas corticis'; $arr[174,1] = 'pest disease'; $arr[174,2] = 'Y'
$arr[175,0] = 'Amynthas tokioensis'; $arr[175,1] = 'pest disease'; $arr[175,2] = 'Y'
$arr[176,0] = 'Amynthas sexpectatus'; $arr[176,1] = 'pest disease'; $arr[176,2] = 'Y'
$arr[177,0] = 'Amynthas gracilis'; $arr[177,1] = 'pest disease'; $arr[177,2] = 'Y'
$arr[178,0] = 'Amynthas diffringens'; $arr[178,1] = 'pest disease'; $arr[178,2] = 'Y'
$arr[179,0] = 'Danaus plexippus'; $arr[179,1] = 'insect'; $arr[179,2] = 'Y'
$arr[180,0] = 'Bombus terricola'; $arr[180,1] = 'insect'; $arr[180,2] = 'Y'
$arr[181,0] = $null; $arr[181,1] = $null; $arr[181,2] = $null
$ws.Range("A1:C182").Value = $arr
Write-Output "Data written"
